# Update "想去人数" (want-to-go count, column F) values on the
# "展览" and "全部类型" worksheets to reflect freshly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    3  = 2704
    4  = 584
    6  = 6632
    7  = 783
    9  = 15
    10 = 12
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
